# Add a new worksheet "Sheet3" at the end of the workbook, populate it with
# the new question/answer data, size its columns, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Header row
$ws3.Range("B1").Value = "What is 40.1 times 3?"
$ws3.Range("C1").Value = "Leeway"

# Data row
$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = 120.3
$ws3.Range("C2").Value = 0.3

# Column widths
$ws3.Columns.Item(2).ColumnWidth = 36.42578125
$ws3.Columns.Item(3).ColumnWidth = 29.5703125

# Selection on the new sheet
[void]$ws3.Range("A3").Select()
